# Update the "APIs: ..." bullet on the Technology list (Subtitle placeholder)
# from "APIs: Cocktail DB, google maps" to a run-split version that reads
# "APIs: Cocktail DB, geoLocation, geoCoding, google maps".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Find the subtitle/content shape that holds the "Technology" bullet list.
$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text.Contains("Cocktail DB")) {
            $targetShape = $shp
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Locate the paragraph that contains the old, single-run text.
$paraCount = $tr.Paragraphs().Count
$targetParaIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $candidate = $tr.Paragraphs($i, 1)
    if ($candidate.Text.Contains("Cocktail DB")) {
        $targetParaIndex = $i
    }
}

$para = $tr.Paragraphs($targetParaIndex, 1)

# The paragraph currently holds a single run: "APIs: Cocktail DB, google maps"
$run1 = $para.Runs(1, 1)

# Rewrite it as 7 runs, matching the target structure:
#   "APIs: Cocktail DB" | ", " | "geoLocation" | ", " | "geoCoding" | ", google " | "maps"
$run1.Text = "APIs: Cocktail DB"
$run2 = $run1.InsertAfter(", ")
$run3 = $run2.InsertAfter("geoLocation")
$run4 = $run3.InsertAfter(", ")
$run5 = $run4.InsertAfter("geoCoding")
$run6 = $run5.InsertAfter(", google ")
$run7 = $run6.InsertAfter("maps")
